# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) totals per game, replacing the previous Strike# counts.
$kValues = @{
    2  = 2
    3  = 3
    4  = 4
    5  = 2
    6  = 5
    7  = 7
    8  = 5
    9  = 2
    10 = 2
    11 = 4
    12 = 7
    13 = 1
    14 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
